$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cliente"
$ws.Range("A2").Value = "Joao"
$ws.Range("A3").Value = "José"
$ws.Range("A4").Value = "Josué"
$ws.Range("B1").Value = "Salario"

$ws.Range("B2").Value = 150
$ws.Range("B3").Value = 180
$ws.Range("B4").Value = 200

$ws.Range("D6").Font.Name = "Calibri"
$ws.Range("D6").Select() | Out-Null
